$d = $word.ActiveDocument

# The "OT (Days)" header cell is built from three runs: "{OT", " (Days)"
# and "}". Remove the " (Days)" run so the cell reads "{OT}".
$rng = $d.Content
[void]$rng.Find.Execute(" (Days)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $editPos = $rng.Start
    $rng.Delete()

    # Word tracks a single hidden "_GoBack" bookmark marking the location
    # of the most recent edit. Re-adding it at the just-edited spot both
    # removes the stale one (previously left at the end of the
    # "Designation:" paragraph) and drops the new, empty one exactly where
    # the deleted text used to be - matching what Word itself does after
    # an in-place edit/save.
    $goBack = $d.Range($editPos, $editPos)
    $d.Bookmarks.Add("_GoBack", $goBack)
}
